$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows for 2023-week (2022-12-16, days 44911) and (2022-12-20, days 44915)
# have been swapped between rows 2-3 and rows 4-5.

# Row 2: previously (44915, Especial, 150, 6000, 6000, 6000, Provincia de Quillota, 1200)
#        now        (44911, Primera,  220, 5000, 5000, 5000, Región de O'Higgins, 1000)
$ws.Range("D2").Value = 44911
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 220
$ws.Range("N2").Value = 5000
$ws.Range("O2").Value = 5000
$ws.Range("P2").Value = 5000
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1000

# Row 3: previously (44915, Primera, 200, 5000, 5000, 5000, Provincia de Quillota, 1000)
#        now        (44911, Segunda, 200, 4000, 4000, 4000, Región de O'Higgins, 800)
$ws.Range("D3").Value = 44911
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 4000
$ws.Range("O3").Value = 4000
$ws.Range("P3").Value = 4000
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 800

# Row 4: previously (44911, Primera, 220, 5000, 5000, 5000, Región de O'Higgins, 1000)
#        now        (44915, Especial, 150, 6000, 6000, 6000, Provincia de Quillota, 1200)
$ws.Range("D4").Value = 44915
$ws.Range("L4").Value = "Especial"
$ws.Range("M4").Value = 150
$ws.Range("N4").Value = 6000
$ws.Range("O4").Value = 6000
$ws.Range("P4").Value = 6000
$ws.Range("R4").Value = "Provincia de Quillota"
$ws.Range("S4").Value = 1200

# Row 5: previously (44911, Segunda, 200, 4000, 4000, 4000, Región de O'Higgins, 800)
#        now        (44915, Primera, 200, 5000, 5000, 5000, Provincia de Quillota, 1000)
$ws.Range("D5").Value = 44915
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 5000
$ws.Range("O5").Value = 5000
$ws.Range("P5").Value = 5000
$ws.Range("R5").Value = "Provincia de Quillota"
$ws.Range("S5").Value = 1000
